$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename header row: "_old" -> "_FV2410", "_new" -> "_FV2504" ---
$newHeaders = @(
    "Segmentname_FV2410","Segmentgruppe_FV2410","Segment_FV2410","Datenelement_FV2410","Segment ID_FV2410",
    "Code_FV2410","Qualifier_FV2410","Beschreibung_FV2410","Bedingungsausdruck_FV2410","Bedingung_FV2410",
    "diff",
    "Segmentname_FV2504","Segmentgruppe_FV2504","Segment_FV2504","Datenelement_FV2504","Segment ID_FV2504",
    "Code_FV2504","Qualifier_FV2504","Beschreibung_FV2504","Bedingungsausdruck_FV2504","Bedingung_FV2504"
)

for ($i = 0; $i -lt $newHeaders.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $newHeaders[$i]
}

# --- Freeze the header row (top row frozen pane) ---
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# --- Turn the used range into a table (ListObject) ---
$rng = $ws.Range("A1:U62")
$lo = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $rng, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$lo.Name = "Table1"

$wb.Save()
